$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.8693910256410257
$ws.Range("C2").Value = 0.9385813148788927
$ws.Range("D2").Value = 0.9026622296173045
$ws.Range("E2").Value = 1156

# Row 3
$ws.Range("B3").Value = 0.9633757961783439
$ws.Range("C3").Value = 0.9322033898305084
$ws.Range("D3").Value = 0.9475332811276428
$ws.Range("E3").Value = 649

# Row 4
$ws.Range("B4").Value = 0.8524173027989822
$ws.Range("C4").Value = 0.850253807106599
$ws.Range("D4").Value = 0.8513341804320204
$ws.Range("E4").Value = 788

# Row 5
$ws.Range("B5").Value = 0.8597122302158273
$ws.Range("C5").Value = 0.6887608069164265
$ws.Range("D5").Value = 0.7648000000000001
$ws.Range("E5").Value = 347

# Row 6 (accuracy row)
$ws.Range("B6").Value = 0.8840136054421769
$ws.Range("C6").Value = 0.8840136054421769
$ws.Range("D6").Value = 0.8840136054421769
$ws.Range("E6").Value = 0.8840136054421769

# Row 7 (macro avg)
$ws.Range("B7").Value = 0.8862240887085447
$ws.Range("C7").Value = 0.8524498296831067
$ws.Range("D7").Value = 0.8665824227942419

# Row 8 (weighted avg)
$ws.Range("B8").Value = 0.8844462230786601
$ws.Range("C8").Value = 0.8840136054421769
$ws.Range("D8").Value = 0.8825386296156041
